$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.064.68"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.679.18"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.70"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.43"
$ws.Range("E9").Value = "  +5.43%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.915.49"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "1.685.27"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.41"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "27.051.39"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.18"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "236.12"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.27"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.52"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.54"
$ws.Range("E27").Value = "  +4.15%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "1.543.90"
$ws.Range("E33").Value = "  +6.67%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  +5.17%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.589"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  +5.27%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.01"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "1.822.13"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.00"
$ws.Range("E50").Value = "  +6.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0508"
$ws.Range("E51").Value = "  +0.21%  "
